# Updates the date line and the 26 division-problem answers in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-03 Tuesday", "2026-02-04 Wednesday"),
    @("69÷9=7, 6", "73÷7=10, 3"),
    @("91÷2=45, 1", "66÷7=9, 3"),
    @("61÷7=8, 5", "37÷6=6, 1"),
    @("68÷9=7, 5", "59÷2=29, 1"),
    @("53÷5=10, 3", "78÷8=9, 6"),
    @("71÷7=10, 1", "29÷4=7, 1"),
    @("32÷5=6, 2", "49÷2=24, 1"),
    @("64÷7=9, 1", "15÷2=7, 1"),
    @("73÷6=12, 1", "92÷6=15, 2"),
    @("20÷7=2, 6", "32÷6=5, 2"),
    @("27÷6=4, 3", "62÷8=7, 6"),
    @("36÷6=6, 0", "24÷6=4, 0"),
    @("42÷7=6, 0", "17÷5=3, 2"),
    @("66÷9=7, 3", "20÷5=4, 0"),
    @("51÷7=7, 2", "53÷4=13, 1"),
    @("14÷4=3, 2", "28÷5=5, 3"),
    @("26÷4=6, 2", "63÷5=12, 3"),
    @("44÷4=11, 0", "12÷9=1, 3"),
    @("76÷3=25, 1", "90÷2=45, 0"),
    @("26÷9=2, 8", "57÷9=6, 3"),
    @("43÷9=4, 7", "42÷5=8, 2"),
    @("46÷7=6, 4", "86÷5=17, 1"),
    @("62÷2=31, 0", "72÷9=8, 0"),
    @("38÷4=9, 2", "70÷7=10, 0"),
    @("48÷4=12, 0", "76÷9=8, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
